# Daily attendance processing - reorder the "Recorded By" (column G) names
# For every row, reverse the order of the comma-separated list of names/emails
# in column G (leaving single-entry cells untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Text

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ',\s*'
        if ($parts.Count -gt 1) {
            $reversed = $parts[($parts.Count - 1)..0]
            $newVal = [string]::Join(', ', $reversed)
            $cell.Value = $newVal
        }
    }
}
